{"js": "// Adds, at the end of the document body, a blank paragraph followed by a\n// paragraph containing a \"Motif du refoulement : \" label and a MERGEFIELD\n// \"=Motif\" merge field (built as a real multi-run Word field: begin /\n// instrText / separate / result / end), matching the formatting already\n// used by the other merge fields in this publipostage template.\n\n// Run-properties shared by every run we create: Times New Roman 12pt,\n// matching the rest of the document's merge-field paragraphs.\nconst rPr =\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:lang w:eastAsia=\"fr-FR\"/>' +\n  '</w:rPr>';\n\n// Same, but with <w:noProof/> added \u2014 used for the field's visible result run.\nconst rPrNoProof =\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:noProof/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:lang w:eastAsia=\"fr-FR\"/>' +\n  '</w:rPr>';\n\nconst pPr = '<w:pPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/>' + rPr + '</w:pPr>';\n\nconst guillemetOpen = String.fromCharCode(0x00ab);\nconst guillemetClose = String.fromCharCode(0x00bb);\n\n// Paragraph 1: empty paragraph (just the paragraph mark formatting).\nconst emptyParagraph = '<w:p>' + pPr + '</w:p>';\n\n// Paragraph 2: label text + the MERGEFIELD field.\nconst motifParagraph =\n  '<w:p>' +\n  pPr +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">Motif du refoulement : </w:t></w:r>' +\n  '<w:r>' + rPr + '<w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r>' + rPr + '<w:instrText xml:space=\"preserve\"> MERGEFIELD  =Motif  \\\\* MERGEFORMAT </w:instrText></w:r>' +\n  '<w:r>' + rPr + '<w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n  '<w:r>' + rPrNoProof + '<w:t>' + guillemetOpen + '=Motif' + guillemetClose + '</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>';\n\nconst bodyFragment =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + emptyParagraph + motifParagraph + '</w:body>' +\n  '</w:document>';\n\n// Office.js requires OOXML passed to insertOoxml to be wrapped in a\n// \"FlatOpc\" package (<pkg:package>...) rather than a bare WordprocessingML\n// fragment.\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' + bodyFragment + '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// Collapse a range to the very end of the document body, then inject both\n// paragraphs there as real WordprocessingML (keeps the exact run/field\n// structure instead of relying on field-insertion heuristics, which build\n// a plain <w:fldSimple>).\nconst endRange = context.document.body.getRange(\"End\");\nendRange.insertOoxml(flatOpc, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Adds, at the end of the document body, a blank paragraph followed by a\n# paragraph containing a \"Motif du refoulement : \" label and a MERGEFIELD\n# \"=Motif\" merge field (built as a real multi-run Word field: begin /\n# instrText / separate / result / end), matching the formatting already\n# used by the other merge fields in this publipostage template.\n\n$d = $word.ActiveDocument\n\n# Run-properties shared by every run we create: Times New Roman 12pt,\n# matching the rest of the document's merge-field paragraphs.\n$rPr = '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:lang w:eastAsia=\"fr-FR\"/>' +\n  '</w:rPr>'\n\n# Same, but with <w:noProof/> added \u2014 used for the field's visible result run.\n$rPrNoProof = '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:noProof/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:lang w:eastAsia=\"fr-FR\"/>' +\n  '</w:rPr>'\n\n$pPr = '<w:pPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/>' + $rPr + '</w:pPr>'\n\n$guillemetOpen = [char]0x00AB\n$guillemetClose = [char]0x00BB\n\n# Paragraph 1: empty paragraph (just the paragraph mark formatting).\n$emptyParagraph = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $pPr + '</w:p>'\n\n# Paragraph 2: label text + the MERGEFIELD field.\n$motifParagraph = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $pPr +\n  '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">Motif du refoulement : </w:t></w:r>' +\n  '<w:r>' + $rPr + '<w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r>' + $rPr + '<w:instrText xml:space=\"preserve\"> MERGEFIELD  =Motif  \\* MERGEFORMAT </w:instrText></w:r>' +\n  '<w:r>' + $rPr + '<w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n  '<w:r>' + $rPrNoProof + '<w:t>' + $guillemetOpen + '=Motif' + $guillemetClose + '</w:t></w:r>' +\n  '<w:r>' + $rPr + '<w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>'\n\n# Collapse a range to the very end of the document body, then inject both\n# paragraphs there as real WordprocessingML (keeps the exact run/field\n# structure instead of relying on the field-insertion heuristics, which\n# build a plain <w:fldSimple>).\n$insertionPoint = $d.Content\n$insertionPoint.Collapse(0)\n$insertionPoint.InsertXML($emptyParagraph + $motifParagraph)\n"}
